$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 header: rename "CBC" column to "CBC1" ---
$ws.Range("F3").Value = "CBC1"

# --- Table 2 headers: split the shared "CBC" column into "CBC2" / "CBC3" ---
$ws.Range("F14").Value = "CBC2"
$ws.Range("G14").Value = "CBC3"

# --- Table 2 body: fill in hours for rows 15-21 (previously blank) ---
$ws.Range("E15:I15").Value = 1.15
$ws.Range("E16:I16").Value = 2.15
$ws.Range("E17:I17").Value = 3.15
$ws.Range("E18:I18").Value = 4.15
$ws.Range("E19:I19").Value = 5.15
$ws.Range("E20:I20").Value = 6.15
$ws.Range("E21:I21").Value = 7.15

# --- selection moves to F5 ---
$ws.Range("F5").Select()
